# Swap the presentation's applied theme color scheme from the custom
# "Integral" (Red Violet) palette to the stock "Office Theme" palette.
#
# (The deck's slide master currently points at a theme part holding the
# Integral/Red Violet clrScheme; the edit replaces those twelve color
# slots with the standard Office theme colors, leaving the font scheme
# and format scheme - already identical between the two themes - untouched.)

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

# Office theme colors (COM RGB = 0xBBGGRR, i.e. byte-reversed hex RRGGBB)
$cs.Item(1).RGB  = 0x000000  # Dark 1      -> 000000
$cs.Item(2).RGB  = 0xFFFFFF  # Light 1     -> FFFFFF
$cs.Item(3).RGB  = 0x6A5444  # Dark 2      -> 44546A
$cs.Item(4).RGB  = 0xE6E6E7  # Light 2     -> E7E6E6
$cs.Item(5).RGB  = 0xD59B5B  # Accent 1    -> 5B9BD5
$cs.Item(6).RGB  = 0x317DED  # Accent 2    -> ED7D31
$cs.Item(7).RGB  = 0xA5A5A5  # Accent 3    -> A5A5A5
$cs.Item(8).RGB  = 0x00C0FF  # Accent 4    -> FFC000
$cs.Item(9).RGB  = 0xC47244  # Accent 5    -> 4472C4
$cs.Item(10).RGB = 0x47AD70  # Accent 6    -> 70AD47
$cs.Item(11).RGB = 0xC16305  # Hyperlink   -> 0563C1
$cs.Item(12).RGB = 0x724F95  # Followed Hyperlink -> 954F72
